$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling (reflects shared-string churn from log -> level transforms) ---
$ws.Range("C1").Value = "y_star_obs Lev"
$ws.Range("D1").Value = "y_star_obs"
$ws.Range("N1").Value = "Y_obs Lev"
$ws.Range("O1").Value = "C_obs Lev"
$ws.Range("P1").Value = "I_obs Lev"
$ws.Range("Q1").Value = "Y_obs"
$ws.Range("R1").Value = "C_obs"
$ws.Range("S1").Value = "I_obs"

# --- Column S (I_obs) data values replaced (level figures -> log-filtered figures) ---
$ws.Range("S2").Value = 0.044572212024319
$ws.Range("S3").Value = 0.019406443418518
$ws.Range("S4").Value = 0.000194719163197
$ws.Range("S5").Value = -0.060480361725084
$ws.Range("S6").Value = -0.035262961891975
$ws.Range("S7").Value = 0.007230245090611
$ws.Range("S8").Value = 0.021472454326682
$ws.Range("S9").Value = -0.017823613695734
$ws.Range("S10").Value = -0.07030171186358
$ws.Range("S11").Value = 0.061078376025932
$ws.Range("S12").Value = 0.008450184330044
$ws.Range("S13").Value = -0.007473095501842
$ws.Range("S14").Value = 0.011480360382841
$ws.Range("S15").Value = 0.013277730237283
$ws.Range("S16").Value = -0.040656091837704
$ws.Range("S17").Value = 0.044259642755579
$ws.Range("S18").Value = 0.027177904879203
$ws.Range("S19").Value = -0.026811891024041
$ws.Range("S20").Value = 0.001717420628818
$ws.Range("S21").Value = 0.0500274104448
$ws.Range("S22").Value = 0.023459125387131
$ws.Range("S23").Value = -0.094896422308285
$ws.Range("S24").Value = -0.125752215215929
$ws.Range("S25").Value = -0.077857762172483
$ws.Range("S26").Value = -0.008081420463805
$ws.Range("S27").Value = -0.031050524805776
$ws.Range("S28").Value = 0.00680376489346
$ws.Range("S29").Value = 0.054375537564084
$ws.Range("S30").Value = 0.135821994722148
$ws.Range("S31").Value = 0.101457679993299
$ws.Range("S32").Value = 0.08902746374612
$ws.Range("S33").Value = 0.068329293683735
$ws.Range("S34").Value = 0.01222799228748
$ws.Range("S35").Value = 0.008111637686184
$ws.Range("S36").Value = -0.016300060933538
$ws.Range("S37").Value = -0.029948846171878
$ws.Range("S38").Value = -0.016511180945443
$ws.Range("S39").Value = -0.024435138828757
$ws.Range("S40").Value = -0.023372154135849
$ws.Range("S41").Value = -0.068606515390746
$ws.Range("S42").Value = -0.057945725313584
$ws.Range("S43").Value = -0.027417736635119
$ws.Range("S44").Value = -0.049044344492506
$ws.Range("S45").Value = -0.003928472988481
$ws.Range("S46").Value = -0.000850425510283
$ws.Range("S47").Value = -0.002428123176438
$ws.Range("S48").Value = 0.027209376164514
$ws.Range("S49").Value = 0.020369812044546
$ws.Range("S50").Value = -0.025693722279525
$ws.Range("S51").Value = -0.002399184615083
$ws.Range("S52").Value = -0.019315792445468
$ws.Range("S53").Value = -0.044404941380257
$ws.Range("S54").Value = 0.005241660528993
$ws.Range("S55").Value = -0.015492956605545
$ws.Range("S56").Value = -0.04034886818653
$ws.Range("S57").Value = 0.038619507609688
$ws.Range("S58").Value = 0.029138860980254
$ws.Range("S59").Value = 0.030199469122348
$ws.Range("S60").Value = 0.042879475041728
$ws.Range("S61").Value = 0.042684090913429
$ws.Range("S62").Value = 0.041478071896492
$ws.Range("S63").Value = 0.047515067371199
$ws.Range("S64").Value = 0.047091634950741
$ws.Range("S65").Value = -0.010022689299904
$ws.Range("S66").Value = -0.020141470023225
$ws.Range("S67").Value = -0.011112156263188
$ws.Range("S68").Value = -0.034108539090559
$ws.Range("S69").Value = 0.015858442340088
$ws.Range("S70").Value = -0.008729186665136
$ws.Range("S71").Value = -0.004361040926929
$ws.Range("S72").Value = 0.018149920368892
$ws.Range("S73").Value = 0.020197268135343
$ws.Range("S74").Value = -0.018640552621775
$ws.Range("S75").Value = 0.03474433902799
$ws.Range("S76").Value = 0.045772038186898
$ws.Range("S77").Value = 0.053842195911383
$ws.Range("S78").Value = 0.072377533105886
$ws.Range("S79").Value = 0.095700761747088
$ws.Range("S80").Value = 0.087405272145897
$ws.Range("S81").Value = 0.034318555953297
$ws.Range("S82").Value = -0.038446609162106
$ws.Range("S83").Value = -0.413942364835563
$ws.Range("S84").Value = -0.177379113915228
$ws.Range("S85").Value = -0.109764777849968
$ws.Range("S86").Value = 0.006102873641478
$ws.Range("S87").Value = -0.0515368999835
$ws.Range("S88").Value = -0.026291322487099
$ws.Range("S89").Value = 0.027112168428213
$ws.Range("S90").Value = 0.100466887394568
$ws.Range("S91").Value = 0.083800209319868
$ws.Range("S92").Value = 0.110037393839118
$ws.Range("S93").Value = 0.111220789903795
$ws.Range("S94").Value = 0.036665898525458
$ws.Range("S95").Value = -0.02823432886607
$ws.Range("S96").Value = -0.048378578761213
$ws.Range("S97").Value = -0.08348042841229
$ws.Range("S98").Value = -0.016922578286943
$ws.Range("S99").Value = -0.002837107599579
$ws.Range("S100").Value = 0.006429852055387
$ws.Range("S101").Value = 0.039697854682615
$ws.Range("S102").Value = -0.013064867371783

# --- Restore the author's last-saved selection on the sheet ---
$ws.Range("S13").Select() | Out-Null
